$wb = $excel.ActiveWorkbook

# --- Insert a new blank worksheet "Sheet1" after "arrayPracticeQuestion" ---
# (ends up positioned before "LinkedList_Try", matching the sheetId=7 / rId5
# shuffle seen in the target workbook.xml)
$afterSheet = $wb.Worksheets.Item("arrayPracticeQuestion")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)

# --- arrayPracticeQuestion: just a selection move (no activation) ---
$arrSheet = $wb.Worksheets.Item("arrayPracticeQuestion")
$arrSheet.Range("D4").Select()

# --- DS-Introduction: new Tree/Graph practice links ---
$dsSheet = $wb.Worksheets.Item("DS-Introduction")

# Add the hyperlink first so its URL text lands in the shared-string table
# before the plain-text URL below (matches shared string ordering 80/81).
$dsSheet.Hyperlinks.Add($dsSheet.Range("D5"), "https://dsportalapp.herokuapp.com/graph/practice")
$dsSheet.Range("D4").Value = "https://dsportalapp.herokuapp.com/Tree/practice"

# Make DS-Introduction the active sheet/tab again, with the updated zoom
# and selection.
$dsSheet.Activate()
$excel.ActiveWindow.Zoom = 78
$dsSheet.Range("D4").Select()

Write-Output "done"
